$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9097349643707275
$ws.Range("B1").Value = 1.744603633880615
$ws.Range("C1").Value = 3.994563579559326
$ws.Range("D1").Value = 3.652950525283813
$ws.Range("E1").Value = 1.507826924324036
